$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds date serial values that were all bumped
# from 45186 (2023-09-17) to 45188 (2023-09-19) for every data row (2..517).
$firstRow = 2
$lastRow = 517

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value = 45188
    }
}
